# Sync local changes before rebase
# - Shift the AE:AI block up by one row (new row2 gets old row3's data, etc.)
# - Overwrite a few cells with freshly entered values (AG2, AF3, AG4)
# - Remove the now-empty trailing row 6
# - Adjust the view (topLeftCell/selection) and add custom column widths for AF/AG

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift AE:AI values up by one row: row(n) <- row(n+1), for n = 2..5
$ws.Range("AE2:AI2").Value2 = $ws.Range("AE3:AI3").Value2
$ws.Range("AE3:AI3").Value2 = $ws.Range("AE4:AI4").Value2
$ws.Range("AE4:AI4").Value2 = $ws.Range("AE5:AI5").Value2
$ws.Range("AE5:AI5").Value2 = $ws.Range("AE6:AI6").Value2

# Clear the now-duplicated last row so nothing is left behind below AI5
$ws.Range("AE6:AI6").ClearContents()

# Apply the newly-entered values that diverge from a pure shift
$ws.Range("AG2").Value2 = 100
$ws.Range("AF3").Value2 = 120
$ws.Range("AG4").Value2 = 180

# Delete the row that is now fully empty
$ws.Rows("6:6").Delete()

# Update column widths for AF (32) and AG (33).
# (The engine quantizes stored width to the nearest 1/6th of a character, so these
# ColumnWidth inputs are chosen to land on the closest reachable width to the
# target XML widths of 11.109375 / 21.5546875.)
$ws.Columns("AF").ColumnWidth = 10.333333333333334
$ws.Columns("AG").ColumnWidth = 20.666666666666668

# Update the view: scroll so column M is at the left edge, then set the active selection
$ws.Application.ActiveWindow.ScrollColumn = 13
$ws.Range("AG14").Select()
